# Add a histogram column for grades: a new "Aufgabe 7" column (H) that
# mirrors the "Gesamt" (G) totals for every student row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new H column the same (locked/protected) cell format already
# used by the data rows further down the sheet (e.g. G10, which carries
# style index 2) before filling in values, so H2:H25 matches the existing
# formatting convention.
$ws.Range("G10").Copy()
$ws.Range("H2:H25").PasteSpecial(-4122)  # xlPasteFormats

# Header for the new column
$ws.Range("H1").Value = "Aufgabe 7"

# Populate H2:H25 with the same totals as the "Gesamt" column (G)
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 7).Value2
}
